$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 72019226
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = ""
$ws.Range("P2").Value = "V om Skarnäs, Upl"
$ws.Range("S2").Value = 5
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2018-05-14"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2018-05-14"
$ws.Range("AF2").NumberFormat = "@"
$ws.Range("AF2").Value = ""
